$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "Admin is redirected to Admin Panel Page"
$ws.Range("C6").Value = "Error message is displayed prompting Admin to fill out the field"
$ws.Range("C8").Value = "An error message is displayed prompting Admin to fill out the empty fields"

$ws.Activate()
$ws.Range("B11").Select()
